# Fix typo in labels: "per capita" -> "per cap." (and bracket/paren fix for
# the "Livestock AB Consumption" label), matching the shared-string updates
# in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old text -> new text, applied to every matching cell on the sheet.
$replacements = @{
    "ln(GDP [dollars per capita])"             = "ln(GDP [dollars per cap.])"
    "ln(Tourism - Inbound [per capita])"       = "ln(Tourism - Inbound [per cap.])"
    "ln(ProMed Mentions [per capita])"         = "ln(ProMed Mentions [per cap.])"
    "ln(Migrant Population [per capita])"      = "ln(Migrant Population [per cap.])"
    "ln(AB Exports [dollars per capita])"      = "ln(AB Exports [dollars per cap.])"
    "ln(Publication Bias Index [per capita])"  = "ln(Publication Bias Index [per cap.])"
    "Livestock AB Consumption [kg per capita)" = "Livestock AB Consumption (kg per cap.)"
}

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($val -ne $null -and $replacements.ContainsKey([string]$val)) {
        $cell.Value2 = $replacements[[string]$val]
    }
}
